$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new header date as plain text (avoid Excel's automatic date parsing)
# via a helper cell + paste-values round trip, so A1 keeps its original style.
$ws.Range("Z1").Value = "'2021/1/15"
$ws.Range("Z1").Copy()
$ws.Range("A1").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").Clear()

# Update the existing attendance entry in row 3 to the new timestamp
$ws.Range("B3").Value = "2021-01-15 09:20:00.484737+00:00"

# Append the rest of the attendance log entries (rows 4-11)
$entries = @(
    @("Muskan Vaswan", "2021-01-15 11:02:48.016808+00:00"),
    @("Muskan Vaswan", "2021-01-19 14:56:54.059903+00:00"),
    @("Muskan Vaswan", "2021-01-19 14:58:27.089588+00:00"),
    @("Muskan Vaswan", "2021-01-19 14:59:06.707655+00:00"),
    @("Muskan Vaswan", "2021-01-19 15:00:43.456282+00:00"),
    @("Muskan Vaswan", "2021-01-19 15:03:03.286919+00:00"),
    @("Muskan Vaswan", "2021-01-19 15:09:10.796135+00:00"),
    @("Muskan Vaswan", "2021-01-20 07:58:50.521689+00:00")
)

$row = 4
foreach ($entry in $entries) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
